# Refresh the crypto price/volume table (columns D and E) with the latest
# scraped snapshot values.
#
# The Price column (D) holds literal text, not numbers (values like
# "29.238.09" or "1.001" are price strings straight from the scrape, several
# of them with more than one "." so they are not even valid numbers). Excel's
# normal cell-input coercion would silently reinterpret a plain-numeric-looking
# replacement (e.g. "0.9998") as a number and drop significant trailing zeros
# (e.g. "240.80" -> 240.8), so for those cells we briefly switch the cell to
# Text format, assign the literal string, then clear the formatting straight
# back off again so the cell ends up holding exactly the source text with no
# left-over cell styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.203.36'
$ws.Cells.Item(2, 5).Value = '  +0.22%  '

$ws.Cells.Item(3, 4).Value = '1.845.28'
$ws.Cells.Item(3, 5).Value = '  +0.70%  '

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9998'
$cell.ClearFormats()
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '240.80'
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -0.18%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6718'
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -1.75%  '

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +0.05%  '

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07420'
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -0.90%  '

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.2956'
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -2.04%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.85'
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -1.31%  '

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07722'
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  +0.82%  '

$ws.Cells.Item(12, 4).Value = '1.814.84'
$ws.Cells.Item(12, 5).Value = '  -0.97%  '

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.013'
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -1.08%  '

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6770'
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -1.04%  '

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '86.15'
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -1.27%  '

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.142'
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -0.38%  '

$ws.Cells.Item(17, 4).Value = '29.197.08'
$ws.Cells.Item(17, 5).Value = '  +0.20%  '

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000008305'
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  +1.39%  '

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '228.80'
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +0.72%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.54'
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -0.12%  '

$ws.Cells.Item(21, 5).Value = '  +0.12%  '

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.200'
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -3.30%  '

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +0.12%  '

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '160.91'
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +0.41%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.686'
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -0.98%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1408'
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -3.59%  '

$ws.Cells.Item(27, 5).Value = '  -0.37%  '

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.509'
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -0.25%  '

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.176'
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  -2.39%  '

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.075'
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -1.77%  '

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.196'
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -0.22%  '

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05335'
$cell.ClearFormats()

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.889'
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +2.34%  '

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7589'
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -1.56%  '

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.138'
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +0.27%  '

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.687'
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  +0.49%  '

$ws.Cells.Item(37, 4).Value = '1.330.38'
$ws.Cells.Item(37, 5).Value = '  +1.23%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.01802'
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -1.75%  '

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.733'
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +0.38%  '

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9250'
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -1.05%  '

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.966'
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +3.02%  '

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +0.25%  '

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '103.47'
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -1.06%  '

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07995'
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +11.91%  '

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.00000000126'
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +2.42%  '

$ws.Cells.Item(46, 4).Value = '1.983.59'
$ws.Cells.Item(46, 5).Value = '  +0.20%  '

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5169'
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -0.51%  '

$ws.Cells.Item(48, 5).Value = '  -0.26%  '

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.11'
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -1.84%  '

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.139'
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -4.75%  '

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05951'
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +0.46%  '
